# Update "want to go" counts (column F) for specific events in the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied identically on both sheets.
$updates = @{
    3  = 1677
    4  = 773
    7  = 11795
    12 = 1106
    14 = 13435
    15 = 13319
    20 = 269
    21 = 90
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
